$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.099.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.78%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.134.88"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.35%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.141.94"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.521"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.151"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.25"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -8.59%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.472"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000252"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -7.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -9.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.643.72"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.52%  "
$ws.Range("E16").Value = "  +1.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.102.16"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.76%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.129.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "479.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.708"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.75"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -7.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -8.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.48"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.12"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -11.92%  "
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.81"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.06%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.114"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -15.63%  "
$ws.Range("B32").Value = "Stacks"
$ws.Range("C32").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.73"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.91%  "
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.25"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.42%  "
$ws.Range("E35").Value = "  -2.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.94"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "53.08"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0754"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.59%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "458.80"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -8.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.95"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -16.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0394"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.50%  "
$ws.Range("E42").Value = "  -9.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.36"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.853.07"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.39%  "
$ws.Range("E45").Value = "  -9.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -11.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "26.20"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -9.87%  "
$ws.Range("E50").Value = "  -4.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "118.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.31%  "
